# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") values are digit-grouped strings (e.g. "25.972.19") that Excel
# would otherwise auto-coerce/round as numbers, so each is written with a temporary
# "@" (text) number format and then restored to the "Normal" style so no stray
# per-cell formatting is left behind. Column B/C/E are plain text already.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '25.972.19'
$r.Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '1.623.05'
$r.Style = "Normal"
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.51%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '213.94'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  -2.44%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.0619'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -3.11%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '18.10'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -7.64%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.0789'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '1.848.14'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '1.644.58'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("E14").Value = '  -2.30%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.522'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -3.87%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '25.948.63'
$r.Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("E17").Value = '  -3.20%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '61.19'
$r.Style = "Normal"
$ws.Range("E18").Value = '  -3.50%  '
$ws.Range("E19").Value = '  +0.52%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '189.49'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -2.73%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '4.23'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("E22").Value = '  -3.73%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '6.06'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -2.40%  '
$ws.Range("E24").Value = '  +0.86%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '143.62'
$r.Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("E28").Value = '  -2.65%  '
$ws.Range("E30").Value = '  -1.71%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '0.0480'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -3.23%  '
$ws.Range("E32").Value = '  -3.84%  '
$ws.Range("E33").Value = '  -5.61%  '
$ws.Range("E34").Value = '  -2.34%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.125.94'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -6.50%  '
$ws.Range("E38").Value = '  -1.26%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.514'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -4.78%  '
$ws.Range("E40").Value = '  -2.35%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '97.63'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -1.62%  '
$ws.Range("E42").Value = '  -2.98%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.758.94'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -1.16%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '5.18'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -5.32%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.0₆0114'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -2.73%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '54.35'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -4.04%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.0524'
$r.Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E50").Value = '  +0.56%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '7.47'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -3.33%  '
